# Re-sort the calibration data rows (2-18) in ascending order of column A (time),
# keeping each row's B/C/D values together with its A value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A2:D18")
$keyRange = $ws.Range("A2:A18")
$dataRange.Sort($keyRange, 1)
